$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H120").Value = 37995
$ws.Range("J120").Value = 37995
$ws.Range("L120").Value = 37995
$ws.Range("N120").Value = -47671
$ws.Range("H127").Value = 1338.1111
$ws.Range("I127").Value = 1218.6
$ws.Range("J127").Value = 1487.5
$ws.Range("K127").Value = 3655.8
$ws.Range("L127").Value = 4462.5
$ws.Range("M127").Value = 1304.2
$ws.Range("N127").Value = -14382.5
$ws.Range("H138").Value = 1603.3835
$ws.Range("I138").Value = 1243.62
$ws.Range("J138").Value = 2385.4783
$ws.Range("K138").Value = 3730.86
$ws.Range("L138").Value = 7156.4349
$ws.Range("M138").Value = 1409.14
$ws.Range("N138").Value = -17436.4349

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6212.37
$ws.Range("I32").Value = 4610.303
$ws.Range("J32").Value = 19174.545
$ws.Range("K32").Value = 4610.303
$ws.Range("L32").Value = 19174.545
$ws.Range("M32").Value = -4323.303
$ws.Range("N32").Value = -19748.545
$ws.Range("H61").Value = 1945.4884
$ws.Range("I61").Value = 2070.3
$ws.Range("J61").Value = 1836.9565
$ws.Range("K61").Value = 2070.3
$ws.Range("L61").Value = 1836.9565
$ws.Range("M61").Value = -1858.3
$ws.Range("N61").Value = -2260.9565
$ws.Range("H122").Value = 1482.7858
$ws.Range("I122").Value = 1296
$ws.Range("J122").Value = 1586.5555
$ws.Range("K122").Value = 3888
$ws.Range("L122").Value = 4759.666499999999
$ws.Range("M122").Value = -1438
$ws.Range("N122").Value = -9659.666499999999
$ws.Range("H136").Value = 1945.4884
$ws.Range("I136").Value = 2070.3
$ws.Range("J136").Value = 1836.9565
$ws.Range("K136").Value = 6210.900000000001
$ws.Range("L136").Value = 5510.8695
$ws.Range("M136").Value = -3660.900000000001
$ws.Range("N136").Value = -10610.8695

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2141.982
$ws.Range("I134").Value = 1254.3143
$ws.Range("J134").Value = 3695.4
$ws.Range("K134").Value = 3762.9429
$ws.Range("L134").Value = 11086.2
$ws.Range("M134").Value = -1227.9429
$ws.Range("N134").Value = -16156.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2944.3684
$ws.Range("I16").Value = 2610.2144
$ws.Range("J16").Value = 3880
$ws.Range("K16").Value = 2610.2144
$ws.Range("L16").Value = 3880
$ws.Range("M16").Value = -2323.2144
$ws.Range("N16").Value = -4454
$ws.Range("H31").Value = 2072.84
$ws.Range("I31").Value = 1131.375
$ws.Range("J31").Value = 3271.068
$ws.Range("K31").Value = 1131.375
$ws.Range("L31").Value = 3271.068
$ws.Range("M31").Value = -836.375
$ws.Range("N31").Value = -3861.068
$ws.Range("H34").Value = 2072.84
$ws.Range("I34").Value = 1131.375
$ws.Range("J34").Value = 3271.068
$ws.Range("K34").Value = 1131.375
$ws.Range("L34").Value = 3271.068
$ws.Range("M34").Value = -929.375
$ws.Range("N34").Value = -3675.068
$ws.Range("H43").Value = 20000
$ws.Range("J43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20368
$ws.Range("H101").Value = 20000
$ws.Range("J101").Value = 20000
$ws.Range("L101").Value = 20000
$ws.Range("N101").Value = -26490
$ws.Range("H113").Value = 2944.3684
$ws.Range("I113").Value = 2610.2144
$ws.Range("J113").Value = 3880
$ws.Range("K113").Value = 2610.2144
$ws.Range("L113").Value = 3880
$ws.Range("M113").Value = -440.2143999999998
$ws.Range("N113").Value = -8220
$ws.Range("H132").Value = 2439.487
$ws.Range("I132").Value = 1713.56
$ws.Range("J132").Value = 3735.7856
$ws.Range("K132").Value = 5140.68
$ws.Range("L132").Value = 11207.3568
$ws.Range("M132").Value = -2610.68
$ws.Range("N132").Value = -16267.3568
$ws.Range("H138").Value = 40380.832
$ws.Range("I138").Value = 19980
$ws.Range("K138").Value = 19980
$ws.Range("M138").Value = -14840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1023.7
$ws.Range("I92").Value = 733.5
$ws.Range("J92").Value = 1217.1666
$ws.Range("K92").Value = 2200.5
$ws.Range("L92").Value = 3651.4998
$ws.Range("M92").Value = -952.5
$ws.Range("N92").Value = -6147.4998
$ws.Range("H107").Value = 50000360
$ws.Range("I107").Value = 71428790
$ws.Range("J107").Value = 674.6667
$ws.Range("K107").Value = 214286370
$ws.Range("L107").Value = 2024.0001
$ws.Range("M107").Value = -214284450
$ws.Range("N107").Value = -5864.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3892.6897
$ws.Range("I16").Value = 2006.0555
$ws.Range("J16").Value = 6979.909
$ws.Range("K16").Value = 2006.0555
$ws.Range("L16").Value = 6979.909
$ws.Range("M16").Value = -1836.0555
$ws.Range("N16").Value = -7319.909
$ws.Range("H40").Value = 29415076
$ws.Range("I40").Value = 47621284
$ws.Range("J40").Value = 5048
$ws.Range("K40").Value = 47621284
$ws.Range("L40").Value = 5048
$ws.Range("M40").Value = -47621148
$ws.Range("N40").Value = -5320
$ws.Range("H46").Value = 588779.7
$ws.Range("I46").Value = 533.8889
$ws.Range("J46").Value = 1250556.2
$ws.Range("K46").Value = 533.8889
$ws.Range("L46").Value = 1250556.2
$ws.Range("M46").Value = -345.8889
$ws.Range("N46").Value = -1250932.2
$ws.Range("H82").Value = 2024.9
$ws.Range("I82").Value = 1407
$ws.Range("J82").Value = 3466.6667
$ws.Range("K82").Value = 1407
$ws.Range("L82").Value = 3466.6667
$ws.Range("M82").Value = -1046
$ws.Range("N82").Value = -4188.6667
$ws.Range("H85").Value = 2024.9
$ws.Range("I85").Value = 1407
$ws.Range("J85").Value = 3466.6667
$ws.Range("K85").Value = 1407
$ws.Range("L85").Value = 3466.6667
$ws.Range("M85").Value = -159
$ws.Range("N85").Value = -5962.6667
$ws.Range("H122").Value = 10145.5
$ws.Range("I122").Value = 35833.332
$ws.Range("J122").Value = 4217.5386
$ws.Range("K122").Value = 107499.996
$ws.Range("L122").Value = 12652.6158
$ws.Range("M122").Value = -105049.996
$ws.Range("N122").Value = -17552.6158

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 27583.334
$ws.Range("J121").Value = 27583.334
$ws.Range("L121").Value = 27583.334
$ws.Range("N121").Value = -31077.334
$ws.Range("H136").Value = 1762.561
$ws.Range("I136").Value = 1658.1482
$ws.Range("J136").Value = 1963.9286
$ws.Range("K136").Value = 4974.444600000001
$ws.Range("L136").Value = 5891.7858
$ws.Range("M136").Value = -2424.444600000001
$ws.Range("N136").Value = -10991.7858
$ws.Range("H141").Value = 56677.5
$ws.Range("J141").Value = 56677.5
$ws.Range("L141").Value = 56677.5
$ws.Range("N141").Value = -67037.5
